$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final row order (row# -> Date, B, C, D) after reordering each
# year block so Oct/Nov/Dec come first, followed by Jan-Sep.
$rows = @(
  @(2, "2014-10", 99.7871, 97.8862, 99.59399999999999),
  @(3, "2014-11", 99.8138, 98.2611, 99.4786),
  @(4, "2014-12", 99.60890000000001, 99.5467, 99.1853),
  @(5, "2014-01", 99.32089999999999, 96.7272, 98.1228),
  @(6, "2014-02", 99.5775, 98.1104, 98.63290000000001),
  @(7, "2014-03", 99.6758, 96.4269, 98.79470000000001),
  @(8, "2014-04", 99.48220000000001, 95.7393, 98.9645),
  @(9, "2014-05", 99.38500000000001, 95.42529999999999, 98.8766),
  @(10, "2014-06", 99.521, 95.4935, 99.0307),
  @(11, "2014-07", 99.65519999999999, 95.71680000000001, 99.3933),
  @(12, "2014-08", 99.65479999999999, 97.0463, 99.40600000000001),
  @(13, "2014-09", 99.8109, 96.76779999999999, 99.4823),
  @(14, "2015-10", 99, 102.9, 98.5),
  @(15, "2015-11", 99.0594, 102.055, 98.43680000000001),
  @(16, "2015-12", 99.2025, 101.1915, 98.52509999999999),
  @(17, "2015-01", 99.518, 99.8489, 98.8528),
  @(18, "2015-02", 99.2384, 98.2671, 98.4689),
  @(19, "2015-03", 98.88330000000001, 97.72450000000001, 98.3823),
  @(20, "2015-04", 99.1014, 98.06319999999999, 98.3069),
  @(21, "2015-05", 99.05200000000001, 99.8596, 98.60290000000001),
  @(22, "2015-06", 98.9873, 100.589, 98.5557),
  @(23, "2015-07", 98.93600000000001, 101.6281, 98.34050000000001),
  @(24, "2015-08", 98.8969, 101.7692, 98.2722),
  @(25, "2015-09", 99.0193, 103.0716, 98.4974),
  @(26, "2016-10", 99.8, 95.40000000000001, 99.5),
  @(27, "2016-11", 100.7, 97.2, 101),
  @(28, "2016-12", 102.4, 98.7, 104.7),
  @(29, "2016-01", 99.4706, 101.7785, 98.6793),
  @(30, "2016-02", 99.60429999999999, 101.7554, 98.92789999999999),
  @(31, "2016-03", 99.64919999999999, 100.919, 99.10509999999999),
  @(32, "2016-04", 99.7715, 98.1575, 98.8896),
  @(33, "2016-05", 99.8, 96.09999999999999, 98.8),
  @(34, "2016-06", 99.8, 94.8, 99),
  @(35, "2016-07", 99.8, 94.40000000000001, 98.90000000000001),
  @(36, "2016-08", 99.59999999999999, 94.59999999999999, 99),
  @(37, "2016-09", 99.7, 94.90000000000001, 99.2),
  @(38, "2017-10", 109.3, 118.1, 121.8),
  @(39, "2017-11", 108.6, 118.6, 119.3),
  @(40, "2017-12", 106.7, 118.7, 113.4),
  @(41, "2017-01", 103.9, 101.7, 107),
  @(42, "2017-02", 104.6, 104.4, 108.2),
  @(43, "2017-03", 104.6, 105.6, 108.3),
  @(44, "2017-04", 104, 109.1, 107.5),
  @(45, "2017-05", 104.1, 110.8, 108.2),
  @(46, "2017-06", 104.6, 110.6, 109.5),
  @(47, "2017-07", 104.9, 110.4, 110.9),
  @(48, "2017-08", 105.8, 108.8, 113.4),
  @(49, "2017-09", 107.5, 114.9, 118.4)
)

foreach ($row in $rows) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
}
